# [FEATURE] Bloqueo y Desbloqueo
# Se agregaron dos casos al modulo y se ajusto configuracion a la base de datos

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet from "Hoja1" to "Database"
$ws.Name = "Database"

# Adjust the database connection configuration:
#  - dBUser value (row 4, column B)
#  - dBPass value (row 5, column B)
$ws.Range("B4").Value = "po3v5snd2tli5v86ntwo"
$ws.Range("B5").Value = "pscale_pw_WDxyGqANy4q29eTUirqYuwo29kwRzaEJFWuO6f0uKqP"

# Reset the active cell selection back to the top of the sheet
$ws.Range("A1").Select()
